# Apply updated values to Sheet1 per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value  = 0.8
$ws.Range("F6").Value  = 0.8272727272727273
$ws.Range("G14").Value = 0.7097744360902256
$ws.Range("G15").Value = 0.6947368421052631
$ws.Range("G16").Value = 0.6947368421052631
$ws.Range("F17").Value = 0.5669172932330827
$ws.Range("G17").Value = 0.6962406015037594
$ws.Range("G18").Value = 0.9055555555555556
$ws.Range("F19").Value = 0.8044444444444444
$ws.Range("G21").Value = 0.8722499999999996
$ws.Range("F22").Value = 0.8934782608695653
$ws.Range("G23").Value = 0.9630434782608696
$ws.Range("F24").Value = 0.8716847826086956
$ws.Range("G24").Value = 0.9652173913043478
$ws.Range("F25").Value = 0.8695652173913043
